$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric; force them to remain text
# to match the original inline-string cell type, by temporarily
# applying a text number format, then clearing formatting again.
$textForceAddrs = @(
    "D5",
    "D6",
    "D10",
    "D13",
    "D17",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D27",
    "D28",
    "D31",
    "D32",
    "D35",
    "D37",
    "D40",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51",
)
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new values
$ws.Range("D2").Value = "43.060.91"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "2.305.45"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "300.67"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "97.92"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  +3.91%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +2.12%  "
$ws.Range("D10").Value = "35.58"
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "17.85"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "2.666.17"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "2.297.17"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "0.787"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "42.977.37"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "13.42"
$ws.Range("E19").Value = "  +8.28%  "
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").Value = "6.13"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").Value = "68.27"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Value = "239.56"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("D27").Value = "24.69"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").Value = "168.35"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  -13.16%  "
$ws.Range("D31").Value = "33.28"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").Value = "5.18"
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("D35").Value = "18.15"
$ws.Range("E35").Value = "  +5.71%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "0.0691"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").Value = "0.111"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").Value = "2.009.54"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").Value = "10.11"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "17.37"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").Value = "54.48"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "2.530.91"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "1.54"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "73.53"
$ws.Range("E51").Value = "  +5.68%  "

# Clear the temporary formatting so these cells end up with
# default (unstyled) formatting, same as the other data cells.
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).ClearFormats()
}
